$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix latitude/longitude data-entry typos for "Center/Centre A Ben Mansour" (id 10001) ---
# G3/G4 had incorrect latitude values (35.x / 36.x) that should match G2 (34.x)
$ws.Range("G3").Value = 34.521169999999998
$ws.Range("G4").Value = 34.521169999999998

# H3/H4 stored near-duplicate longitude text values; normalize them to the same
# text value as H2 so the redundant shared strings can be dropped.
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = " -6.453275"
$ws.Range("H3").Style = "Normal"

$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = " -6.453275"
$ws.Range("H4").Style = "Normal"

# --- Update number_of_kiosks (column L) ---
# id 10001 (rows 2-4): 3 kiosks
$ws.Range("L2:L4").Value = 3

# id 10002-10015 (rows 5-46): 2 kiosks
$ws.Range("L5:L46").Value = 2

# Restore the editor selection left behind by the authoring session
$ws.Range("L28").Select() | Out-Null
